$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 7-11, mirroring the layout of row 6 (columns A-D blank, E-L/N identical
# static values), but with their own Email (col M) and Name (col O) values.
$newRows = @(
    @{ Row = 7;  Email = "testuser11@gmail.com"; Name = "Uday" },
    @{ Row = 8;  Email = "testuser12@gmail.com"; Name = "Rohit" },
    @{ Row = 9;  Email = "testuser13@gmail.com"; Name = "Cleo" },
    @{ Row = 10; Email = "testuser14@gmail.com"; Name = "Monty" },
    @{ Row = 11; Email = "testuser15@gmail.com"; Name = "Emilly" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy formatting + values from the template row (row 6) into the new row.
    $ws.Range("A6:O6").Copy($ws.Range("A" + $r + ":O" + $r))

    # Match the row height used by the existing data rows.
    $ws.Rows.Item($r).RowHeight = 18

    # Overwrite the per-row values: Email (col M) and Name (col O).
    $ws.Cells.Item($r, 13).Value = $item.Email
    $ws.Cells.Item($r, 15).Value = $item.Name
}

$excel.CutCopyMode = 0
